# Update countries & provincias Spain
# Refresh of the COVID dashboard data (new pull at 04:46 replacing the 03:29 pull),
# plus a reordering of two country-name pairs in the shared-strings table:
#   - "Paises Bajos" / "Belgica" swap rank (rows 29/30)
#   - "Montserrat" / "Islas Malvinas" swap rank (rows 216/217)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Octubre de 2020 a las 04:46"

# --- Row 27: Pakistan ---------------------------------------------------
$ws.Range("B27").Value = 327063
$ws.Range("C27").Value = 847
$ws.Range("D27").Value = 310101
$ws.Range("E27").Value = 10235
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 12
$ws.Range("H27").Value = 6727

# --- Rows 29/30: Paises Bajos & Belgica swap rank, new figures --------
$ws.Range("A29").Value = "Belgica"
$ws.Range("B29").Value = 287700
$ws.Range("C29").Value = 17568
$ws.Range("D29").Value = 22517
$ws.Range("E29").Value = 254525
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 70
$ws.Range("H29").Value = 10658

$ws.Range("A30").Value = "Paises Bajos"
$ws.Range("B30").Value = 272401
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 6964

# --- Row 38: Bolivia -----------------------------------------------------
$ws.Range("B38").Value = 140612
$ws.Range("C38").Value = 167
$ws.Range("D38").Value = 107633
$ws.Range("E38").Value = 24371
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 24
$ws.Range("H38").Value = 8608

# --- Row 190: Camboya ------------------------------------------------
$ws.Range("B190").Value = 287
$ws.Range("C190").Value = 1
$ws.Range("D190").Value = 283
$ws.Range("E190").Value = 4
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 0

# --- Rows 216/217: Montserrat & Islas Malvinas swap rank ---------------
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 13
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0

$ws.Range("A217").Value = "Montserrat"
$ws.Range("B217").Value = 13
$ws.Range("C217").Value = 0
$ws.Range("D217").Value = 12
$ws.Range("E217").Value = 0
$ws.Range("F217").Value = 0
$ws.Range("G217").Value = 0
$ws.Range("H217").Value = 1
